$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) from the existing header cell H1 onto the
# two new header cells so they pick up the same bold/border/alignment style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data columns I2:I40 and J2:J40 ---
$iValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,9,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,6,1)
$jValues = @(4,4,5,6,7,5,7,6,3,6,5,5,6,7,9,6,6,6,5,6,6,6,5,6,5,6,5,5,6,4,4,4,4,5,4,5,3,7,2)

for ($i = 0; $i -lt 39; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
